$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "[논문 요약] On the Role of Bidirectionality in Language Model Pre-Training"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222746609163"

$ws.Range("D27").Value = "알라꿍달라꿍의 대화요약 이모저모"
$ws.Range("E27").Value = "https://blog.pingpong.us/alaggung-dlaggung-dialog-summary/"

$ws.Range("D37").Value = "[Paper Review] Recurrent Reconstructive Network for Sequential Anomaly Detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1984&mod=document&pageid=1"
